function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 3 ("Reference project"): resize the title + content boxes
# and refresh the dataset-size bullet / add a new bullet point.
# ---------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

$titleShape = Get-ShapeById $slide3 2
$titleShape.Height = 84.72890093779507

$contentShape = Get-ShapeById $slide3 7
$contentShape.Top = 163.7871653543307
$contentShape.Height = 336.68370078740156

$tr = $contentShape.TextFrame.TextRange
$datasetPara = $tr.Characters(135, 111)
$datasetPara.Text = "Dataset was taken from UCI website and represented only selected number of observations – around 11k out of 45k."
[void]$tr.InsertAfter("`r The core aim is to predict if the customer will subscribe or not")

# ---------------------------------------------------------------
# Slide 16 ("AI tools"): fix the duplicated "use" typo.
# ---------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$aiShape = Get-ShapeById $slide16 7
$aiTr = $aiShape.TextFrame.TextRange
$aiFirstLine = $aiTr.Characters(1, 58)
$aiFirstLine.Text = "We used following AI tools while building our project:"
